$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new "2020" data column (N), mirroring the formatting of the
# existing last column (M) for each of the three data rows.

$ws.Range("M3").Copy()
$ws.Range("N3").PasteSpecial(-4122)
$ws.Range("N3").Value = 2020

$ws.Range("M4").Copy()
$ws.Range("N4").PasteSpecial(-4122)
$ws.Range("N4").Value = 15

$ws.Range("M5").Copy()
$ws.Range("N5").PasteSpecial(-4122)
$ws.Range("N5").Value = 1308.3

$excel.CutCopyMode = $false

# Move the active selection to the newly added cell, matching the
# workbook's last saved selection state.
$ws.Range("N6").Select()
